$d = $word.ActiveDocument
$d.Content.Find.Execute("957÷6=159, 3", $true, $false, $false, $false, $false, $true, 1, $false, "679÷6=113, 1", 2) | Out-Null
$d.Content.Find.Execute("904÷4=226, 0", $true, $false, $false, $false, $false, $true, 1, $false, "217÷2=108, 1", 2) | Out-Null
$d.Content.Find.Execute("710÷8=88, 6", $true, $false, $false, $false, $false, $true, 1, $false, "425÷6=70, 5", 2) | Out-Null
$d.Content.Find.Execute("491÷3=163, 2", $true, $false, $false, $false, $false, $true, 1, $false, "398÷7=56, 6", 2) | Out-Null
$d.Content.Find.Execute("181÷7=25, 6", $true, $false, $false, $false, $false, $true, 1, $false, "591÷2=295, 1", 2) | Out-Null
$d.Content.Find.Execute("584÷5=116, 4", $true, $false, $false, $false, $false, $true, 1, $false, "639÷7=91, 2", 2) | Out-Null
$d.Content.Find.Execute("108÷9=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "957÷5=191, 2", 2) | Out-Null
$d.Content.Find.Execute("759÷9=84, 3", $true, $false, $false, $false, $false, $true, 1, $false, "572÷8=71, 4", 2) | Out-Null
$d.Content.Find.Execute("694÷2=347, 0", $true, $false, $false, $false, $false, $true, 1, $false, "484÷8=60, 4", 2) | Out-Null
$d.Content.Find.Execute("748÷8=93, 4", $true, $false, $false, $false, $false, $true, 1, $false, "276÷8=34, 4", 2) | Out-Null
$d.Content.Find.Execute("876÷5=175, 1", $true, $false, $false, $false, $false, $true, 1, $false, "215÷2=107, 1", 2) | Out-Null
$d.Content.Find.Execute("417÷7=59, 4", $true, $false, $false, $false, $false, $true, 1, $false, "499÷9=55, 4", 2) | Out-Null
$d.Content.Find.Execute("788÷2=394, 0", $true, $false, $false, $false, $false, $true, 1, $false, "702÷8=87, 6", 2) | Out-Null
$d.Content.Find.Execute("601÷4=150, 1", $true, $false, $false, $false, $false, $true, 1, $false, "269÷7=38, 3", 2) | Out-Null
$d.Content.Find.Execute("479÷2=239, 1", $true, $false, $false, $false, $false, $true, 1, $false, "347÷6=57, 5", 2) | Out-Null
$d.Content.Find.Execute("612÷4=153, 0", $true, $false, $false, $false, $false, $true, 1, $false, "522÷4=130, 2", 2) | Out-Null
$d.Content.Find.Execute("517÷3=172, 1", $true, $false, $false, $false, $false, $true, 1, $false, "433÷9=48, 1", 2) | Out-Null
$d.Content.Find.Execute("407÷4=101, 3", $true, $false, $false, $false, $false, $true, 1, $false, "913÷6=152, 1", 2) | Out-Null
$d.Content.Find.Execute("237÷6=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "225÷8=28, 1", 2) | Out-Null
$d.Content.Find.Execute("338÷3=112, 2", $true, $false, $false, $false, $false, $true, 1, $false, "355÷9=39, 4", 2) | Out-Null
$d.Content.Find.Execute("464÷4=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "301÷7=43, 0", 2) | Out-Null
$d.Content.Find.Execute("306÷3=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "973÷3=324, 1", 2) | Out-Null
$d.Content.Find.Execute("621÷7=88, 5", $true, $false, $false, $false, $false, $true, 1, $false, "769÷4=192, 1", 2) | Out-Null
$d.Content.Find.Execute("340÷9=37, 7", $true, $false, $false, $false, $false, $true, 1, $false, "138÷4=34, 2", 2) | Out-Null
$d.Content.Find.Execute("156÷3=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "787÷6=131, 1", 2) | Out-Null
